$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.632.66"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.208.86"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.99%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  -0.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.203.56"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.128"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.66"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.402"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.754.76"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("E14").Value = "  +1.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.49"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.550.44"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000167"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.199.73"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("E19").Value = "  +1.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.82%  "

$ws.Range("E23").Value = "  -0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.62"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.27%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000118"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.502"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.328.48"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.74"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.10%  "

$ws.Range("E29").Value = "  +1.46%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.94"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.51"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.63%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.29"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.70"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.21"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "168.35"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.27%  "

$ws.Range("E38").Value = "  +3.50%  "

$ws.Range("E39").Value = "  +2.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.39"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.56"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.00%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.30"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.04%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.668.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.23"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0668"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "328.09"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.39%  "

$ws.Range("E50").Value = "  +2.56%  "

$ws.Range("E51").Value = "  +0.35%  "
